$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 17725.545
$ws.Range("J18").Value = 17996.8
$ws.Range("L18").Value = 17996.8
$ws.Range("N18").Value = -18564.8
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H74").Value = 5250
$ws.Range("I74").Value = 5333.3335
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 5333.3335
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -4397.3335
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 5250
$ws.Range("I77").Value = 5333.3335
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 26666.6675
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -21986.6675
$ws.Range("N77").Value = -34360
$ws.Range("H132").Value = 772.26086
$ws.Range("I132").Value = 696.4286
$ws.Range("K132").Value = 2089.2858
$ws.Range("M132").Value = 440.7142000000003
$ws.Range("H138").Value = 1625.9596
$ws.Range("I138").Value = 1009
$ws.Range("J138").Value = 2062.0862
$ws.Range("K138").Value = 3027
$ws.Range("L138").Value = 6186.258600000001
$ws.Range("M138").Value = 2113
$ws.Range("N138").Value = -16466.2586
$ws.Range("H141").Value = 4002892.2
$ws.Range("J141").Value = 5851.6665
$ws.Range("L141").Value = 17554.9995
$ws.Range("N141").Value = -27914.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19488.8
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 19488.8
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 19488.8
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -20062.8
$ws.Range("H74").Value = 1037.2333
$ws.Range("I74").Value = 811.6923
$ws.Range("K74").Value = 811.6923
$ws.Range("M74").Value = 62.30769999999995
$ws.Range("H77").Value = 1037.2333
$ws.Range("I77").Value = 811.6923
$ws.Range("K77").Value = 4058.4615
$ws.Range("M77").Value = 309.5384999999997
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 1174.0667
$ws.Range("I122").Value = 1174.0667
$ws.Range("K122").Value = 3522.2001
$ws.Range("M122").Value = -1072.2001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2680
$ws.Range("I20").Value = 2650.3215
$ws.Range("J20").Value = 2743.923
$ws.Range("K20").Value = 2650.3215
$ws.Range("L20").Value = 2743.923
$ws.Range("M20").Value = -2403.3215
$ws.Range("N20").Value = -3237.923
$ws.Range("H94").Value = 1021.3889
$ws.Range("I94").Value = 492.46667
$ws.Range("J94").Value = 3666
$ws.Range("K94").Value = 492.46667
$ws.Range("L94").Value = 3666
$ws.Range("M94").Value = -41.46667000000002
$ws.Range("N94").Value = -4568
$ws.Range("H134").Value = 5988.5757
$ws.Range("I134").Value = 7647.591
$ws.Range("K134").Value = 22942.773
$ws.Range("M134").Value = -20407.773

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1001
$ws.Range("I7").Value = 1001
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1001
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -888
$ws.Range("N7").ClearContents()
$ws.Range("H31").Value = 2304.0527
$ws.Range("I31").Value = 1677.1
$ws.Range("J31").Value = 3000.6667
$ws.Range("K31").Value = 1677.1
$ws.Range("L31").Value = 3000.6667
$ws.Range("M31").Value = -1382.1
$ws.Range("N31").Value = -3590.6667
$ws.Range("H34").Value = 2304.0527
$ws.Range("I34").Value = 1677.1
$ws.Range("J34").Value = 3000.6667
$ws.Range("K34").Value = 1677.1
$ws.Range("L34").Value = 3000.6667
$ws.Range("M34").Value = -1475.1
$ws.Range("N34").Value = -3404.6667
$ws.Range("H59").Value = 18142.857
$ws.Range("J59").Value = 18142.857
$ws.Range("L59").Value = 18142.857
$ws.Range("N59").Value = -20432.857
$ws.Range("H105").Value = 1084.1538
$ws.Range("I105").Value = 1130.8
$ws.Range("J105").Value = 928.6667
$ws.Range("K105").Value = 1130.8
$ws.Range("L105").Value = 928.6667
$ws.Range("M105").Value = 616.2
$ws.Range("N105").Value = -4422.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 558.96
$ws.Range("I5").Value = 540.3333
$ws.Range("J5").Value = 606.8570999999999
$ws.Range("K5").Value = 1620.9999
$ws.Range("L5").Value = 1820.5713
$ws.Range("M5").Value = -1508.9999
$ws.Range("N5").Value = -2044.5713
$ws.Range("H81").Value = 3000.8333
$ws.Range("J81").Value = 3000.8333
$ws.Range("L81").Value = 9002.499899999999
$ws.Range("N81").Value = -11248.4999
$ws.Range("H84").Value = 3000.8333
$ws.Range("J84").Value = 3000.8333
$ws.Range("L84").Value = 27007.4997
$ws.Range("N84").Value = -38239.4997
$ws.Range("H107").Value = 427.57144
$ws.Range("I107").Value = 363.8889
$ws.Range("J107").Value = 475.33334
$ws.Range("K107").Value = 1091.6667
$ws.Range("L107").Value = 1426.00002
$ws.Range("M107").Value = 828.3333
$ws.Range("N107").Value = -5266.000019999999
$ws.Range("H118").Value = 2114.7
$ws.Range("J118").Value = 3198.6
$ws.Range("L118").Value = 9595.799999999999
$ws.Range("N118").Value = -12081.8
$ws.Range("H121").Value = 663.2308
$ws.Range("J121").Value = 737.55554
$ws.Range("L121").Value = 2212.66662
$ws.Range("N121").Value = -4832.66662
$ws.Range("H131").Value = 15096.286
$ws.Range("J131").Value = 15831.6
$ws.Range("L131").Value = 47494.8
$ws.Range("N131").Value = -57574.8
$ws.Range("H135").Value = 558.96
$ws.Range("I135").Value = 540.3333
$ws.Range("J135").Value = 606.8570999999999
$ws.Range("K135").Value = 4862.9997
$ws.Range("L135").Value = 5461.7139
$ws.Range("M135").Value = -2327.9997
$ws.Range("N135").Value = -10531.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5618001
$ws.Range("J12").Value = 1817503
$ws.Range("L12").Value = 1817503
$ws.Range("N12").Value = -1817783
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50518
$ws.Range("H97").Value = 1064.3846
$ws.Range("I97").Value = 1074.8636
$ws.Range("K97").Value = 1074.8636
$ws.Range("M97").Value = -578.8635999999999
$ws.Range("H102").Value = 3249.25
$ws.Range("I102").Value = 5498.5
$ws.Range("K102").Value = 5498.5
$ws.Range("M102").Value = -3876.5
$ws.Range("H113").Value = 1997
$ws.Range("J113").Value = 1997
$ws.Range("L113").Value = 1997
$ws.Range("N113").Value = -6337
$ws.Range("H122").Value = 1375.4706
$ws.Range("J122").Value = 1470.7142
$ws.Range("L122").Value = 4412.142599999999
$ws.Range("N122").Value = -9312.142599999999
$ws.Range("H132").Value = 5499308.5
$ws.Range("I132").Value = 7696114
$ws.Range("J132").Value = 7293.5
$ws.Range("K132").Value = 23088342
$ws.Range("L132").Value = 21880.5
$ws.Range("M132").Value = -23085812
$ws.Range("N132").Value = -26940.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1197.5333
$ws.Range("I22").Value = 487.8
$ws.Range("J22").Value = 2617
$ws.Range("K22").Value = 487.8
$ws.Range("L22").Value = 2617
$ws.Range("M22").Value = -192.8
$ws.Range("N22").Value = -3207
$ws.Range("H27").Value = 1197.5333
$ws.Range("I27").Value = 487.8
$ws.Range("J27").Value = 2617
$ws.Range("K27").Value = 487.8
$ws.Range("L27").Value = 2617
$ws.Range("M27").Value = -380.8
$ws.Range("N27").Value = -2831
$ws.Range("H40").Value = 10292.111
$ws.Range("I40").Value = 10433.5
$ws.Range("J40").Value = 9797.25
$ws.Range("K40").Value = 10433.5
$ws.Range("L40").Value = 9797.25
$ws.Range("M40").Value = -10297.5
$ws.Range("N40").Value = -10069.25
$ws.Range("H46").Value = 2832.3635
$ws.Range("I46").Value = 1799.6666
$ws.Range("J46").Value = 3219.625
$ws.Range("K46").Value = 1799.6666
$ws.Range("L46").Value = 3219.625
$ws.Range("M46").Value = -1611.6666
$ws.Range("N46").Value = -3595.625
$ws.Range("H61").Value = 1663.6
$ws.Range("I61").Value = 1617.7273
$ws.Range("K61").Value = 1617.7273
$ws.Range("M61").Value = -1415.7273
$ws.Range("H113").Value = 1663.6
$ws.Range("I113").Value = 1617.7273
$ws.Range("K113").Value = 1617.7273
$ws.Range("M113").Value = 552.2727
$ws.Range("H132").Value = 2074.238
$ws.Range("I132").Value = 1786.35
$ws.Range("K132").Value = 5359.049999999999
$ws.Range("M132").Value = -2829.049999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1155.25
$ws.Range("I100").Value = 774.75
$ws.Range("J100").Value = 1535.75
$ws.Range("K100").Value = 1549.5
$ws.Range("L100").Value = 3071.5
$ws.Range("M100").Value = -1008.5
$ws.Range("N100").Value = -4153.5
$ws.Range("H108").Value = 66999.5
$ws.Range("J108").Value = 66999.5
$ws.Range("L108").Value = 66999.5
$ws.Range("N108").Value = -74679.5
$ws.Range("H126").Value = 2125.4443
$ws.Range("I126").Value = 1911.6428
$ws.Range("K126").Value = 5734.928400000001
$ws.Range("M126").Value = -3264.928400000001
